# Add "Washing Produce" content to the bottom of the sheet, then insert a new
# "food_produce" / "Produce" row near the top of the TCS food section.
#
# The order of writes below matters: it reproduces the exact order the
# original author typed the content in (title row, then key column, then
# value column, then numbers, and only at the very end inserting the new
# row up in the TCS section) so that new shared-string entries land on the
# same indices as the authoritative edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "Washing Produce" table, appended right after the existing data
#        (old last row was 46, so the new content starts at row 47). ---

# Title row
$ws.Range("A47").Value = "washProduce_title"
$ws.Range("B47").Value = "Washing Produce"
$ws.Range("C47").Value = 2.5

# Item keys (column A) for the remaining 7 rows
$ws.Range("A48").Value = "washProduce_item_1"
$ws.Range("A49").Value = "washProduce_item_2"
$ws.Range("A50").Value = "washProduce_item_3"
$ws.Range("A51").Value = "washProduce_item_4"
$ws.Range("A52").Value = "washProduce_item_5"
$ws.Range("A53").Value = "washProduce_item_6"
$ws.Range("A54").Value = "washProduce_item_7"

# Item text (column B) for the remaining 7 rows
$ws.Range("B48").Value = "Wash your hands for 20 seconds with warm water and soap before and after preparing fresh produce."
$ws.Range("B49").Value = "If damage or bruising occurs before eating or handling, cut away the damaged or bruised areas before preparing or eating."
$ws.Range("B50").Value = "Rinse produce BEFORE you peel it, so dirt and bacteria aren’t transferred from the knife onto the fruit or vegetable."
$ws.Range("B51").Value = "Gently rub produce while holding under plain running water. There’s no need to use soap or a produce wash."
$ws.Range("B52").Value = "Use a clean vegetable brush to scrub firm produce, such as melons and cucumbers."
$ws.Range("B53").Value = "Dry produce with a clean cloth or paper towel to further reduce bacteria that may be present."
$ws.Range("B54").Value = "Remove the outermost leaves of a head of lettuce or cabbage."

# Durations (column C) for the remaining 7 rows
$ws.Range("C48").Value = 4
$ws.Range("C49").Value = 6
$ws.Range("C50").Value = 6
$ws.Range("C51").Value = 6
$ws.Range("C52").Value = 6
$ws.Range("C53").Value = 6
$ws.Range("C54").Value = 5

# --- 2. Insert a new "Produce" row into the TCS foods list at row 29,
#        pushing everything from the old row 29 onward (including the
#        table we just typed) down by one row. ---

$ws.Rows("29:29").Insert()
$ws.Range("A29").Value = "food_produce"
$ws.Range("B29").Value = "Produce"

# --- 3. Update the view's active selection to match the author's final
#        cursor position. ---

$ws.Range("A56").Select()
